$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal TEXT (matches the source file, where
# every data cell -- including numeric-looking prices/percentages -- is
# stored as an inline string, not a Number). Forcing NumberFormat='@' before
# the assignment stops Excel from auto-coercing strings like "236.66" into
# a float, and resetting the Style back to "Normal" afterwards drops the
# temporary text format so no extra cell formatting is left behind.
function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Updated cryptocurrency market data (price + 1h volume change). Rows 31/32
# and 47/48 were re-ranked, so the coin name/link/price/volume that used to
# sit in one row now belongs to the other.

# Row 2
Set-TextValue 2 4 "92.354.29"
Set-TextValue 2 5 "  +0.67%  "

# Row 3
Set-TextValue 3 4 "3.093.19"
Set-TextValue 3 5 "  -2.23%  "

# Row 4
Set-TextValue 4 5 "  -0.01%  "

# Row 5
Set-TextValue 5 4 "236.66"
Set-TextValue 5 5 "  -1.28%  "

# Row 6
Set-TextValue 6 4 "609.83"
Set-TextValue 6 5 "  -1.97%  "

# Row 7
Set-TextValue 7 4 "1.09"
Set-TextValue 7 5 "  -3.36%  "

# Row 8
Set-TextValue 8 4 "0.389"
Set-TextValue 8 5 "  +3.32%  "

# Row 9
Set-TextValue 9 5 "  -0.05%  "

# Row 10
Set-TextValue 10 4 "3.089.95"
Set-TextValue 10 5 "  -2.34%  "

# Row 11
Set-TextValue 11 4 "0.734"
Set-TextValue 11 5 "  -1.50%  "

# Row 12
Set-TextValue 12 5 "  -1.80%  "

# Row 13
Set-TextValue 13 4 "0.0000247"
Set-TextValue 13 5 "  -0.16%  "

# Row 14
Set-TextValue 14 4 "92.279.30"
Set-TextValue 14 5 "  +1.15%  "

# Row 15
Set-TextValue 15 4 "34.08"
Set-TextValue 15 5 "  -4.31%  "

# Row 16
Set-TextValue 16 4 "5.43"
Set-TextValue 16 5 "  -2.96%  "

# Row 17
Set-TextValue 17 4 "3.668.32"
Set-TextValue 17 5 "  -2.28%  "

# Row 18
Set-TextValue 18 4 "3.104.81"
Set-TextValue 18 5 "  -2.50%  "

# Row 19
Set-TextValue 19 4 "3.76"
Set-TextValue 19 5 "  -0.19%  "

# Row 20
Set-TextValue 20 4 "14.59"
Set-TextValue 20 5 "  -4.92%  "

# Row 21
Set-TextValue 21 4 "5.73"
Set-TextValue 21 5 "  -3.71%  "

# Row 22
Set-TextValue 22 4 "9.31"
Set-TextValue 22 5 "  +0.94%  "

# Row 23
Set-TextValue 23 4 "443.66"
Set-TextValue 23 5 "  -2.92%  "

# Row 24
Set-TextValue 24 4 "0.0000195"
Set-TextValue 24 5 "  -5.15%  "

# Row 25
Set-TextValue 25 4 "5.69"
Set-TextValue 25 5 "  -5.90%  "

# Row 26
Set-TextValue 26 4 "85.83"
Set-TextValue 26 5 "  -3.73%  "

# Row 27
Set-TextValue 27 4 "11.63"
Set-TextValue 27 5 "  -3.80%  "

# Row 28
Set-TextValue 28 4 "3.254.43"
Set-TextValue 28 5 "  -1.93%  "

# Row 30
Set-TextValue 30 4 "0.131"
Set-TextValue 30 5 "  +3.66%  "

# Row 31
Set-TextValue 31 2 "Cronos"
Set-TextValue 31 3 "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue 31 4 "0.169"
Set-TextValue 31 5 "  -2.14%  "

# Row 32
Set-TextValue 32 2 "Stellar"
Set-TextValue 32 3 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue 32 4 "0.226"
Set-TextValue 32 5 "  -2.39%  "

# Row 33
Set-TextValue 33 4 "9.10"
Set-TextValue 33 5 "  -3.55%  "

# Row 34
Set-TextValue 34 5 "  +6.60%  "

# Row 35
Set-TextValue 35 4 "7.87"
Set-TextValue 35 5 "  +2.09%  "

# Row 36
Set-TextValue 36 5 "  -7.64%  "

# Row 37
Set-TextValue 37 4 "25.87"
Set-TextValue 37 5 "  -2.83%  "

# Row 38
Set-TextValue 38 4 "3.87"
Set-TextValue 38 5 "  +1.06%  "

# Row 39
Set-TextValue 39 4 "1.89"
Set-TextValue 39 5 "  -3.54%  "

# Row 40
Set-TextValue 40 4 "482.91"
Set-TextValue 40 5 "  -6.47%  "

# Row 41
Set-TextValue 41 4 "23.92"
Set-TextValue 41 5 "  +7.76%  "

# Row 42
Set-TextValue 42 4 "1.28"
Set-TextValue 42 5 "  -5.92%  "

# Row 43
Set-TextValue 43 4 "0.429"
Set-TextValue 43 5 "  -5.60%  "

# Row 44
Set-TextValue 44 4 "3.31"
Set-TextValue 44 5 "  -4.92%  "

# Row 46
Set-TextValue 46 4 "163.58"
Set-TextValue 46 5 "  +3.06%  "

# Row 47
Set-TextValue 47 2 "Stacks"
Set-TextValue 47 3 "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue 47 4 "1.86"
Set-TextValue 47 5 "  -4.33%  "

# Row 48
Set-TextValue 48 2 "ARBITRUM"
Set-TextValue 48 3 "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue 48 4 "0.682"
Set-TextValue 48 5 "  -4.75%  "

# Row 49
Set-TextValue 49 4 "1.37"
Set-TextValue 49 5 "  -0.74%  "

# Row 50
Set-TextValue 50 4 "0.0329"
Set-TextValue 50 5 "  +2.04%  "

# Row 51
Set-TextValue 51 4 "43.93"
Set-TextValue 51 5 "  -0.38%  "
